# Natmi following Dr Hou advice
# Rebuild the Lif-Il6st LR-pair sheet with the full 3x3 sending/target cluster cross-product
# (FAPs, sCs, ECs) instead of the original 3-row same-index-only pairing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("FAPs","Lif","Il6st","FAPs",1,0.3333333333333333,0.5355756666666667,1.606727,0.168705706571243,0.168705706571243,3,1,39.96608766666667,119.898263,0.2616165719423124,0.2616165719423124,21.40486404613345,192.643776415201,0.04413620862027424,0.04413620862027425),
    @("FAPs","Lif","Il6st","sCs",1,0.3333333333333333,0.5355756666666667,1.606727,0.168705706571243,0.168705706571243,3,1,85.119611,255.358833,0.5571899111219771,0.557189911121977,45.58799240773234,410.291931669591,0.09400111765020125,0.09400111765020124),
    @("FAPs","Lif","Il6st","ECs",1,0.3333333333333333,0.5355756666666667,1.606727,0.168705706571243,0.168705706571243,3,1,27.68018833333333,83.040565,0.1811935169357105,0.1811935169357105,14.82483532008389,133.423517880755,0.03056838030076752,0.03056838030076752),
    @("sCs","Lif","Il6st","FAPs",3,1,2.335066666666667,7.0052,0.7355432601013561,0.7355432601013562,3,1,39.96608766666667,119.898263,0.2616165719423124,0.2616165719423124,93.32347910751112,839.9113119676001,0.1924303062229894,0.1924303062229895),
    @("sCs","Lif","Il6st","sCs",3,1,2.335066666666667,7.0052,0.7355432601013561,0.7355432601013562,3,1,85.119611,255.358833,0.5571899111219771,0.557189911121977,198.7599663257333,1788.8396969316,0.4098372837222439,0.4098372837222439),
    @("sCs","Lif","Il6st","ECs",3,1,2.335066666666667,7.0052,0.7355432601013561,0.7355432601013562,3,1,27.68018833333333,83.040565,0.1811935169357105,0.1811935169357105,64.63508510422223,581.715765938,0.1332756701561227,0.1332756701561228),
    @("ECs","Lif","Il6st","FAPs",3,1,0.3039726666666667,0.911918,0.09575103332740086,0.09575103332740086,3,1,39.96608766666667,119.898263,0.2616165719423124,0.2616165719423124,12.14859824427045,109.337384198434,0.02505005709904872,0.02505005709904872),
    @("ECs","Lif","Il6st","sCs",3,1,0.3039726666666667,0.911918,0.09575103332740086,0.09575103332740086,3,1,85.119611,255.358833,0.5571899111219771,0.557189911121977,25.87403514129933,232.866316271694,0.05335150974953196,0.05335150974953194),
    @("ECs","Lif","Il6st","ECs",3,1,0.3039726666666667,0.911918,0.09575103332740086,0.09575103332740086,3,1,27.68018833333333,83.040565,0.1811935169357105,0.1811935169357105,8.414020661518888,75.72618595367,0.01734946647882018,0.01734946647882018)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowNum = $i + 2
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowData[$c]
    }
}